$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WorkSheet 1")

# New row of data to append as row 19, matching the style/format of the row above it
$row = 19

$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 42601.899050925924
$ws.Cells.Item($row, 2).Value = "Random"
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 48
$ws.Cells.Item($row, 9).Value = 52
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 65
$ws.Cells.Item($row, 13).Value = 35
